$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.772.22'
$ws.Range("E2").Value = '  -1.63%  '
$ws.Range("D3").Value = '3.695.01'
$ws.Range("E3").Value = '  -2.32%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '615.43'
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.61'
$ws.Range("E6").Value = '  -0.81%  '
$ws.Range("D7").Value = '3.693.60'
$ws.Range("E7").Value = '  -2.35%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("E9").Value = '  -2.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.164'
$ws.Range("E10").Value = '  -2.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.24'
$ws.Range("E11").Value = '  -2.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.479'
$ws.Range("E12").Value = '  -4.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.72'
$ws.Range("E13").Value = '  -3.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000253'
$ws.Range("E14").Value = '  -2.30%  '
$ws.Range("D15").Value = '4.306.44'
$ws.Range("E15").Value = '  -2.43%  '
$ws.Range("D16").Value = '3.695.98'
$ws.Range("E16").Value = '  -2.51%  '
$ws.Range("D17").Value = '69.723.39'
$ws.Range("E17").Value = '  -1.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.122'
$ws.Range("E18").Value = '  -1.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.53'
$ws.Range("E19").Value = '  -1.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.35'
$ws.Range("E20").Value = '  -2.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '500.87'
$ws.Range("E21").Value = '  -4.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.13'
$ws.Range("E22").Value = '  -3.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.712'
$ws.Range("E23").Value = '  -5.04%  '
$ws.Range("E24").Value = '  +1.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.23'
$ws.Range("E25").Value = '  -2.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.39'
$ws.Range("E26").Value = '  +3.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.93'
$ws.Range("E27").Value = '  -5.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000129'
$ws.Range("E28").Value = '  +4.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.43'
$ws.Range("E30").Value = '  -3.44%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.88'
$ws.Range("E31").Value = '  -1.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.95'
$ws.Range("E32").Value = '  -1.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.06'
$ws.Range("E33").Value = '  -7.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.114'
$ws.Range("E34").Value = '  -2.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.03'
$ws.Range("E37").Value = '  -2.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.137'
$ws.Range("E38").Value = '  +3.07%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.337'
$ws.Range("E39").Value = '  -1.58%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.94'
$ws.Range("E40").Value = '  -2.84%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.05'
$ws.Range("E41").Value = '  -7.99%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.93'
$ws.Range("E42").Value = '  +4.26%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '429.38'
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("B44").Value = 'Arweave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '44.37'
$ws.Range("E44").Value = '  +0.74%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.54'
$ws.Range("E45").Value = '  -4.00%  '
$ws.Range("D46").Value = '2.935.65'
$ws.Range("E46").Value = '  -7.10%  '
$ws.Range("E47").Value = '  -2.91%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.23'
$ws.Range("E48").Value = '  -2.37%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.78'
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.42'
$ws.Range("E51").Value = '  -2.70%  '
